$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph that contains a unique piece of text and
# return its 1-based Paragraphs index (robust to earlier edits shifting
# character offsets).
# ---------------------------------------------------------------------------
function Get-ParaIndex($text) {
    $rng = $d.Content
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $p = $rng.Paragraphs.Item(1)
    return $p.Index
}

# ---------------------------------------------------------------------------
# 1. Remove the blank paragraph that currently separates item "22." from
#    item "23.". Deleting the blank paragraph's Range (it holds nothing but
#    its own paragraph mark) merges the two neighbouring paragraphs'
#    boundaries exactly like the target: "22." is immediately followed by
#    "23." with no blank line between them.
# ---------------------------------------------------------------------------
$idx22 = Get-ParaIndex("22. Moved DisplayLabyrinth method to the DisplayMaze method in Maze class. LabyrinthGenerator method move to GenerateMaze method in MazeClass")
$blankAfter22 = $d.Paragraphs.Item($idx22 + 1)
$blankAfter22.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Item "31.": insert the missing space between "rendering" and
#    "(drawing)" so the sentence reads "...responsible for rendering
#    (drawing) on the console."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("rendering(drawing)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "rendering (drawing)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Item "32.": replace the whole sentence with the new wording.
# ---------------------------------------------------------------------------
$idx32 = Get-ParaIndex("32. Moved all game constants in the introduced structure GameConstants. Removed all magic strings and numbers and defined as constants in the GameConstants structure.")
$p32 = $d.Paragraphs.Item($idx32)
$fullRange = $d.Range($p32.Range.Start, $p32.Range.End - 1)
$fullRange.Text = "32. Replaced all magic strings and numbers and defined as constants in the corresponding class."

# ---------------------------------------------------------------------------
# 4. Insert a brand new item "33." paragraph right after item "32.".
# ---------------------------------------------------------------------------
$idx32 = Get-ParaIndex("32. Replaced all magic strings and numbers and defined as constants in the corresponding class.")
$p32 = $d.Paragraphs.Item($idx32)
$p32.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx32 + 1)
$newPara.Range.Text = "33. Renamed variable s in AddScore method in the ScoreBoard class to currentPlayer and renamed variables s1 and s2 in SortScores method in the ScoreBoard class respectively to currentPlayer and otherPlayer."

# ---------------------------------------------------------------------------
# 5. Relocate the "_GoBack" bookmark so that it wraps the end of the "22."
#    paragraph (instead of living in its own trailing paragraph). Adding a
#    bookmark named "_GoBack" moves the existing singleton bookmark rather
#    than creating a duplicate, leaving its old paragraph perfectly empty -
#    which becomes one of the two blank paragraphs at the end of the
#    document, matching the target structure.
# ---------------------------------------------------------------------------
$idx22 = Get-ParaIndex("22. Moved DisplayLabyrinth method to the DisplayMaze method in Maze class. LabyrinthGenerator method move to GenerateMaze method in MazeClass")
$p22 = $d.Paragraphs.Item($idx22)
$bmRange = $d.Range($p22.Range.Start, $p22.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
